$wb = $excel.ActiveWorkbook

# Overview sheet: Latest HO Xliff Generate Date
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-26 15:19:18"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-26 15:19:13"
$wsZhCn.Range("K2").Value = "2016-08-26 15:19:30"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-26 15:19:18"
$wsDeDe.Range("K2").Value = "2016-08-26 15:19:37"
